# FY23 Points Tracker Audit - Dec update, Q3 start
# Fill in actuals for the Q3 section (rows 26/27) and a few engagement
# "x" achievement marks (rows 46/51/52). All subtotal/total rows
# (24, 44, 58, 61) are driven by existing formulas and recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Q3 "Revenue" row (26) ---
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 15
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 4
$ws.Range("I26").Value = 4

# --- Q3 "Logos" row (27) ---
$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 10
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 5
$ws.Range("I27").Value = 5

# --- Eng achievement marks ---
$ws.Range("I46").Value = "x"
$ws.Range("I51").Value = "x"
$ws.Range("D52").Value = "x"

# --- PTO Bonus (20) row (56) ---
$ws.Range("C56").Value = 24
$ws.Range("D56").Value = 34
$ws.Range("E56").Value = 24
$ws.Range("F56").Value = 21
$ws.Range("G56").Value = 20
$ws.Range("I56").Value = 16

# --- Kudos (5 / 1) row (57) ---
$ws.Range("C57").Value = 3
$ws.Range("D57").Value = 6
$ws.Range("H57").Value = 5
$ws.Range("I57").Value = 9

# Scroll/select so the frozen pane + active cell match the saved view
# (Excel recomputes topLeftCell for the visible window automatically).
$ws.Range("G46").Select()
